$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 72381256
$ws.Range("I61").Value = 7143229
$ws.Range("J61").Value = 333333340
$ws.Range("K61").Value = 21429687
$ws.Range("L61").Value = 1000000020
$ws.Range("M61").Value = -21429515
$ws.Range("N61").Value = -1000000364
$ws.Range("H70").Value = 2936.25
$ws.Range("I70").Value = 1996
$ws.Range("K70").Value = 5988
$ws.Range("M70").Value = -5718
$ws.Range("H73").Value = 2936.25
$ws.Range("I73").Value = 1996
$ws.Range("K73").Value = 5988
$ws.Range("M73").Value = -5052
$ws.Range("H129").Value = 4263.25
$ws.Range("J129").Value = 4027.1538
$ws.Range("L129").Value = 12081.4614
$ws.Range("N129").Value = -22081.4614
$ws.Range("H131").Value = 5547.6665
$ws.Range("I131").Value = 911.1667
$ws.Range("J131").Value = 10184.167
$ws.Range("K131").Value = 2733.5001
$ws.Range("L131").Value = 30552.501
$ws.Range("M131").Value = 2306.4999
$ws.Range("N131").Value = -40632.501
$ws.Range("H135").Value = 33338522
$ws.Range("I135").Value = 33338522
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 300046698
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -300044163
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 2981.0603
$ws.Range("I138").Value = 1422.9
$ws.Range("K138").Value = 4268.700000000001
$ws.Range("M138").Value = 871.2999999999993

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 51249.832
$ws.Range("I32").Value = 36489.566
$ws.Range("K32").Value = 36489.566
$ws.Range("M32").Value = -36202.566
$ws.Range("H61").Value = 3949.5557
$ws.Range("I61").Value = 3877.4167
$ws.Range("J61").Value = 4526.6665
$ws.Range("K61").Value = 3877.4167
$ws.Range("L61").Value = 4526.6665
$ws.Range("M61").Value = -3665.4167
$ws.Range("N61").Value = -4950.6665
$ws.Range("H74").Value = 83353090
$ws.Range("J74").Value = 166701330
$ws.Range("L74").Value = 166701330
$ws.Range("N74").Value = -166703078
$ws.Range("H77").Value = 83353090
$ws.Range("J77").Value = 166701330
$ws.Range("L77").Value = 833506650
$ws.Range("N77").Value = -833515386
$ws.Range("H122").Value = 103761.1
$ws.Range("I122").Value = 146767.28
$ws.Range("J122").Value = 3413.3333
$ws.Range("K122").Value = 440301.84
$ws.Range("L122").Value = 10239.9999
$ws.Range("M122").Value = -437851.84
$ws.Range("N122").Value = -15139.9999
$ws.Range("H136").Value = 3949.5557
$ws.Range("I136").Value = 3877.4167
$ws.Range("J136").Value = 4526.6665
$ws.Range("K136").Value = 11632.2501
$ws.Range("L136").Value = 13579.9995
$ws.Range("M136").Value = -9082.250100000001
$ws.Range("N136").Value = -18679.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 600683
$ws.Range("J43").Value = 600683
$ws.Range("L43").Value = 600683
$ws.Range("N43").Value = -601045
$ws.Range("H86").Value = 6889.615
$ws.Range("I86").Value = 3446.375
$ws.Range("K86").Value = 3446.375
$ws.Range("M86").Value = -2323.375
$ws.Range("H89").Value = 6889.615
$ws.Range("I89").Value = 3446.375
$ws.Range("K89").Value = 17231.875
$ws.Range("M89").Value = -11615.875
$ws.Range("H105").Value = 7287.8423
$ws.Range("I105").Value = 6968.7646
$ws.Range("K105").Value = 6968.7646
$ws.Range("M105").Value = -5221.7646
$ws.Range("H107").Value = 1683.3334
$ws.Range("I107").Value = 1990.7
$ws.Range("J107").Value = 1299.125
$ws.Range("K107").Value = 1990.7
$ws.Range("L107").Value = 1299.125
$ws.Range("M107").Value = -70.70000000000005
$ws.Range("N107").Value = -5139.125
$ws.Range("H134").Value = 2883.5
$ws.Range("I134").Value = 2725.4614
$ws.Range("J134").Value = 3568.3333
$ws.Range("K134").Value = 8176.3842
$ws.Range("L134").Value = 10704.9999
$ws.Range("M134").Value = -5641.3842
$ws.Range("N134").Value = -15774.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6609799
$ws.Range("I31").Value = 3015923
$ws.Range("K31").Value = 3015923
$ws.Range("M31").Value = -3015628
$ws.Range("H33").Value = 19060.143
$ws.Range("I33").Value = 1357.75
$ws.Range("J33").Value = 42663.332
$ws.Range("K33").Value = 1357.75
$ws.Range("L33").Value = 42663.332
$ws.Range("M33").Value = -978.75
$ws.Range("N33").Value = -43421.332
$ws.Range("H34").Value = 6609799
$ws.Range("I34").Value = 3015923
$ws.Range("K34").Value = 3015923
$ws.Range("M34").Value = -3015721
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H86").Value = 506548.06
$ws.Range("I86").Value = 594744.75
$ws.Range("K86").Value = 594744.75
$ws.Range("M86").Value = -593621.75
$ws.Range("H89").Value = 506548.06
$ws.Range("I89").Value = 594744.75
$ws.Range("K89").Value = 2973723.75
$ws.Range("M89").Value = -2968107.75
$ws.Range("H134").Value = 15758.796
$ws.Range("I134").Value = 15758.796
$ws.Range("K134").Value = 47276.388
$ws.Range("M134").Value = -44741.388

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2430.1428
$ws.Range("J5").Value = 7700
$ws.Range("L5").Value = 23100
$ws.Range("N5").Value = -23324
$ws.Range("H33").Value = 3029.6924
$ws.Range("I33").Value = 799.5
$ws.Range("J33").Value = 4941.2856
$ws.Range("K33").Value = 4797
$ws.Range("L33").Value = 29647.7136
$ws.Range("M33").Value = -4514
$ws.Range("N33").Value = -30213.7136
$ws.Range("H92").Value = 749.5
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2997
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5493
$ws.Range("H97").Value = 146.16667
$ws.Range("I97").Value = 117.5
$ws.Range("J97").Value = 160.5
$ws.Range("K97").Value = 352.5
$ws.Range("L97").Value = 481.5
$ws.Range("M97").Value = 143.5
$ws.Range("N97").Value = -1473.5
$ws.Range("H113").Value = 385
$ws.Range("J113").Value = 270
$ws.Range("L113").Value = 810
$ws.Range("N113").Value = -5150
$ws.Range("H117").Value = 980.4
$ws.Range("I117").Value = 472.2857
$ws.Range("J117").Value = 2166
$ws.Range("K117").Value = 1416.8571
$ws.Range("L117").Value = 6498
$ws.Range("M117").Value = 2025.1429
$ws.Range("N117").Value = -13382
$ws.Range("H131").Value = 37397.355
$ws.Range("I131").Value = 84202.586
$ws.Range("J131").Value = 2293.4375
$ws.Range("K131").Value = 252607.758
$ws.Range("L131").Value = 6880.3125
$ws.Range("M131").Value = -247567.758
$ws.Range("N131").Value = -16960.3125
$ws.Range("H135").Value = 2430.1428
$ws.Range("J135").Value = 7700
$ws.Range("L135").Value = 69300
$ws.Range("N135").Value = -74370

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142859650
$ws.Range("I80").Value = 333335260
$ws.Range("J80").Value = 2944
$ws.Range("K80").Value = 333335260
$ws.Range("L80").Value = 2944
$ws.Range("M80").Value = -333334262
$ws.Range("N80").Value = -4940
$ws.Range("H83").Value = 142859650
$ws.Range("I83").Value = 333335260
$ws.Range("J83").Value = 2944
$ws.Range("K83").Value = 1666676300
$ws.Range("L83").Value = 14720
$ws.Range("M83").Value = -1666671308
$ws.Range("N83").Value = -24704
$ws.Range("H97").Value = 2169.389
$ws.Range("I97").Value = 2091.1177
$ws.Range("J97").Value = 3500
$ws.Range("K97").Value = 2091.1177
$ws.Range("L97").Value = 3500
$ws.Range("M97").Value = -1595.1177
$ws.Range("N97").Value = -4492

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3412.6667
$ws.Range("I16").Value = 3654.2727
$ws.Range("K16").Value = 3654.2727
$ws.Range("M16").Value = -3484.2727
$ws.Range("H35").Value = 3545.8572
$ws.Range("I35").Value = 2304.2
$ws.Range("J35").Value = 6650
$ws.Range("K35").Value = 2304.2
$ws.Range("L35").Value = 6650
$ws.Range("M35").Value = -1968.2
$ws.Range("N35").Value = -7322
$ws.Range("H68").Value = 2625.8667
$ws.Range("I68").Value = 2149.077
$ws.Range("J68").Value = 5725
$ws.Range("K68").Value = 2149.077
$ws.Range("L68").Value = 5725
$ws.Range("M68").Value = -1400.077
$ws.Range("N68").Value = -7223
$ws.Range("H71").Value = 2625.8667
$ws.Range("I71").Value = 2149.077
$ws.Range("J71").Value = 5725
$ws.Range("K71").Value = 10745.385
$ws.Range("L71").Value = 28625
$ws.Range("M71").Value = -7001.385000000002
$ws.Range("N71").Value = -36113
$ws.Range("H122").Value = 6372.8335
$ws.Range("I122").Value = 5059.25
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 15177.75
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -12727.75
$ws.Range("N122").Value = -31900

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17818.438
$ws.Range("I62").Value = 17857
$ws.Range("J62").Value = 17548.5
$ws.Range("K62").Value = 17857
$ws.Range("L62").Value = 17548.5
$ws.Range("M62").Value = -17233
$ws.Range("N62").Value = -18796.5
$ws.Range("H65").Value = 17818.438
$ws.Range("I65").Value = 17857
$ws.Range("J65").Value = 17548.5
$ws.Range("K65").Value = 89285
$ws.Range("L65").Value = 87742.5
$ws.Range("M65").Value = -86165
$ws.Range("N65").Value = -93982.5
$ws.Range("H96").Value = 2533.3333
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 4600
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 4600
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -7346
$ws.Range("H122").Value = 9049.75
$ws.Range("I122").Value = 4142.7144
$ws.Range("J122").Value = 15919.6
$ws.Range("K122").Value = 12428.1432
$ws.Range("L122").Value = 47758.8
$ws.Range("M122").Value = -9978.143199999999
$ws.Range("N122").Value = -52658.8
$ws.Range("H136").Value = 1893.0416
$ws.Range("I136").Value = 1655.9474
$ws.Range("J136").Value = 2794
$ws.Range("K136").Value = 4967.8422
$ws.Range("L136").Value = 8382
$ws.Range("M136").Value = -2417.8422
$ws.Range("N136").Value = -13482
